$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B width changes from 15.42578125 to 14.7109375 (OOXML "width" units)
$ws.Columns.Item(2).ColumnWidth = 13.833333333333334

# Update data values in columns A and B, rows 1-32
$ws.Cells.Item(1, 1).Value = -0.23420996532603056
$ws.Cells.Item(1, 2).Value = 0.23371492428795193
$ws.Cells.Item(2, 1).Value = -0.18760660231405168
$ws.Cells.Item(2, 2).Value = 0.18571353667604029
$ws.Cells.Item(3, 1).Value = -0.082769337886027472
$ws.Cells.Item(3, 2).Value = 0.082325110598400641
$ws.Cells.Item(4, 1).Value = -0.14431760423700268
$ws.Cells.Item(4, 2).Value = 0.14355324445096684
$ws.Cells.Item(5, 1).Value = -0.13755324525437018
$ws.Cells.Item(5, 2).Value = 0.13602560049804335
$ws.Cells.Item(6, 1).Value = -0.063188142736221309
$ws.Cells.Item(6, 2).Value = 0.06312885919851885
$ws.Cells.Item(7, 1).Value = -0.043128860164300065
$ws.Cells.Item(7, 2).Value = 0.04302700470288201
$ws.Cells.Item(8, 1).Value = -0.023027005674100209
$ws.Cells.Item(8, 2).Value = 0.02298651174180133
$ws.Cells.Item(9, 1).Value = -0.047323867009231968
$ws.Cells.Item(9, 2).Value = 0.047156676477919213
$ws.Cells.Item(10, 1).Value = 0.013856519110838406
$ws.Cells.Item(10, 2).Value = -0.013856270041010532
$ws.Cells.Item(11, 1).Value = 0.018356269207043141
$ws.Cells.Item(11, 2).Value = -0.018373926355060632
$ws.Cells.Item(12, 1).Value = -0.045399291701738509
$ws.Cells.Item(12, 2).Value = 0.045151957799666409
$ws.Cells.Item(13, 1).Value = -0.03915195866118637
$ws.Cells.Item(13, 2).Value = 0.039084901965933838
$ws.Cells.Item(14, 1).Value = -0.027084902889598084
$ws.Cells.Item(14, 2).Value = 0.027052856656778701
$ws.Cells.Item(15, 1).Value = -0.021052857524884949
$ws.Cells.Item(15, 2).Value = 0.021027665859398859
$ws.Cells.Item(16, 1).Value = -0.01502766673009992
$ws.Cells.Item(16, 2).Value = 0.015004150292793561
$ws.Cells.Item(17, 1).Value = -0.0090041511670486685
$ws.Cells.Item(17, 2).Value = 0.0089999990959492848
$ws.Cells.Item(18, 1).Value = -0.036108322786059688
$ws.Cells.Item(18, 2).Value = 0.036096323006379549
$ws.Cells.Item(19, 1).Value = -0.027096323809120193
$ws.Cells.Item(19, 2).Value = 0.027013237414548641
$ws.Cells.Item(20, 1).Value = -0.018013238224833472
$ws.Cells.Item(20, 2).Value = 0.018004247913053106
$ws.Cells.Item(21, 1).Value = -0.0090042487243744418
$ws.Cells.Item(21, 2).Value = 0.0089999991880418406
$ws.Cells.Item(22, 1).Value = -0.093944199608277046
$ws.Cells.Item(22, 2).Value = 0.093632759074433025
$ws.Cells.Item(23, 1).Value = -0.084632759899911036
$ws.Cells.Item(23, 2).Value = 0.084126288165847996
$ws.Cells.Item(24, 1).Value = -0.042126289326080979
$ws.Cells.Item(24, 2).Value = 0.041999998833591512
$ws.Cells.Item(25, 1).Value = -0.066837458589272813
$ws.Cells.Item(25, 2).Value = 0.066721913548459355
$ws.Cells.Item(26, 1).Value = -0.060721914377907638
$ws.Cells.Item(26, 2).Value = 0.060576350039337257
$ws.Cells.Item(27, 1).Value = -0.054576350871679224
$ws.Cells.Item(27, 2).Value = 0.054087734977851021
$ws.Cells.Item(28, 1).Value = -0.048087735821471078
$ws.Cells.Item(28, 2).Value = 0.047763883671461116
$ws.Cells.Item(29, 1).Value = -0.035763884580763516
$ws.Cells.Item(29, 2).Value = 0.035623775383536227
$ws.Cells.Item(30, 1).Value = -0.015623776373540732
$ws.Cells.Item(30, 2).Value = 0.01560439997383245
$ws.Cells.Item(31, 1).Value = -0.063415048352279513
$ws.Cells.Item(31, 2).Value = 0.063279579565856281
$ws.Cells.Item(32, 1).Value = -0.0060005105841902662
$ws.Cells.Item(32, 2).Value = 0.0059999991399246611
